# Generate Report for Handback
#
# Updates the timestamp columns on the "Overview", "zh-cn" and "de-de"
# sheets of the handback-status report to reflect the latest
# handoff/handback generation times.

$wb = $excel.ActiveWorkbook

# Overview sheet, row 2 ("f9db31a2-...md"): "Latest HO Xliff Generate Date"
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-23 11:06:39"

# zh-cn sheet, row 2: "Correspond Handoff Datetime" / "Correspond Handback DateTime"
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-23 11:06:34"
$wsZhCn.Range("K2").Value = "2016-08-23 11:06:51"

# de-de sheet, row 2: "Correspond Handoff Datetime" / "Correspond Handback DateTime"
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-23 11:06:39"
$wsDeDe.Range("K2").Value = "2016-08-23 11:06:58"
